$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Apps Data"

$ws.Range("A1").Value = "App Name"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Downloads"
$ws.Range("D1").Value = "Developer Email"
$ws.Range("E1").Value = "Category"
$ws.Range("F1").Value = "App Link"
$ws.Range("G1").Value = "Website"

$ws.Range("A2").Value = "Piere: Budget & Manage Money"
$ws.Range("B2").Value = "Discover intelligent money management and budgeting that works. Achieve your goals with individually tailored financial insights, custom budgets, smarter shared payments, and spend tracking. Quickly view all your accounts, monthly bills, create goals, and send money while building knowledge of your finances.`nPiere is the budgeting, personal finance and money management app that becomes a single place to manage all your finances. View account balances and track budgets while setting financial goals, and as Piere learns from your personal finances, it automatically adjusts your budget to achieve your goals. Piere is quickly becoming the trusted budgeting app that more users depend upon to consistently deliver intelligent money management and spending insights.`nA SINGLE APP FOR ALL YOUR MONEY`nBring all your banking, spending, and financial accounts into a single place. Piere ensures you have a more accurate view of your financial health by merging all of your: budgets, monthly expenses, account balances, spending, net worth, and more. Connect your checking and savings, cash, credit cards, loans, investments, mortgages, retirement savings, and more.`nUNDERSTAND YOUR CASH FLOW`nPiere’s intelligence tracks your budgets, transactions, monthly expenses, and subscriptions. The app calculates your spending trends and net worth from a more complete picture of your personal finances. You’ll be notified when you’re about to go over budget or overdraft an account, or when a recurring subscription price increases. Piere’s intelligence then automatically adjusts your budget to ensure you still meet your goals, making managing your budget easy and empowering.`nSPEND AND SAVE MORE, CONFIDENTLY WITH PIERE INSIGHTS`nPiere’s rapidly evolving intelligence analyzes your budgets, spending habits, and accounts to recommend new ways to create money goals that you can achieve. Use Piere’s money tracker to get a smarter view of how you spend and save.`nSMARTER BUDGETING AND EXPENSE MONITORING`nOur budgeting feature is powered by personalized, evolving financial intelligence. In 2-clicks get a smart budget tailored to your unique spending and never again spend hours to create and manage your budget.`nYOUR BILLS AND EXPENSES IN A NEW LIGHT`nWith all your accounts, transactions, and balances in once place, Piere is a smarter bill tracker and a simpler way to manage your expenses. Keep an eye on your debt and avoid late fees with a single view of all your accounts. Your goals are finally achievable with Piere’s intelligent expense tracking.`nPERSONALIZED FIANANCIAL REPORTS`nEffortlessly customize spending and income reports by time, category, account, tags, and more, with charts that are interactive and intuitive.`nWIN AT YOUR FINANCIAL GOALS`nCreate financial goals such as spending and investing targets and get customized, actionable insights from Piere’s budgeting app. Smarter money management advice from Piere helps you make progress while living a richer life today. Piere is here to help you realize the full potential of your money.`nU.S. & CANADIAN BUDGETING APP`nIn addition to covering over 98% of U.S. financial institutions, Piere is made with Canadian personal finance management in mind. Connect a rapidly growing number of Canadian banks and financial institutions and take advantage of monthly budgets, expense tracking, net worth monitoring and more to empower your financial management.`nPRIVATE AND SECURE`nOur users' data security is top of mind for us. We keep your data safe and your account protected by encrypting your data at rest and in transit, and completed the same security and compliance process as the big banks. Review our full privacy and security policy here: https://piere.com/privacypolicy"
$ws.Range("C2").Value = "1,000+"
$ws.Range("D2").Value = "contact@piere.com"
$ws.Range("E2").Value = "fintech"
$ws.Range("F2").Value = "https://play.google.com/store/apps/details?id=com.pieremobile"
$ws.Range("G2").Value = "https://piere.com"

$ws.Range("A3").Value = "Bilance - Your Money & Budget"
$ws.Range("B3").Value = "Make your personal finances make sense with Bilance: the most automated money manager.`nBilance connects securely with your bank account and automatically categorizes your expenses with a powerful AI model. Take advantage of easy budgeting features to take control of your spending and build better financial habits.`nWhether you want to reduce financial stress, save more or start investing, you’ll love keeping track of your money with Bilance. Get started today with a 30-day free trial! `nWhy Bilance?`n• Bilance saves you time and gives you peace of mind. Automatic bank sync and AI-powered expense categorization keep you up to date with your cash flow – no manual work needed.`n• Bilance makes your finances make sense. Leave the math for us – you can simply review your spending trends on clear infographics and automatic reports.`n• Bring all your finances together. Connect as many banks and bank accounts with Bilance as you want, including neobanks and cash.`n• Focused on your security. We don’t sell your data and we have built Bilance using a secure infrastructure.`n• Works across Europe. Over 800 banks in 31 countries already support automatic sync.`nMore key features:`n• Compare with others - see how your expenses stack up`n• Savings rate - monitor your savings rate in real-time`n• Budgets - set monthly budgets and control your spending`n• Recurring - tracker for your subscriptions and bills`n• Reports - get an automatic monthly report at the end of each month`n• Tags - organize your spending`n• Dark mode - exquisite looks, also at nighttime`n• No ads - we don’t show any ads or recommend other products`n• Quick support - contact us anytime through chat`n• Constantly improving - we release new features and improvements regularly`nWhat users are saying about Bilance:`n“One of the greatest money management apps I’ve seen so far. Bank account sync works really neat and makes tracking convenient.” – Andres`n“It’s amazing that the app automatically sets a category for all your transactions.” – Ingrid`n“I have used countless other personal finance apps, but Bilance has everything! Definitely recommend it!” – Marta`n“Thanks for making me think about my money.” – Herbert`nSubscription information:`n- Try Bilance for free for 30 days. After that, you can choose between a Monthly and Yearly subscription.`n- Payments will be charged to your Apple account`n- The subscription is auto-renewable`n- Your account will be charged for renewal within 24 hours before the end of your current subscription period.`n- Auto-renewal may be turned off at any moment in the App Store settings.`nLearn more: https://www.bilanceapp.com`nContact us: hello@bilanceapp.com`nTerms of Use: https://www.bilanceapp.com/terms`nPrivacy Policy: https://www.bilanceapp.com/privacy"
$ws.Range("C3").Value = "10,000+"
$ws.Range("D3").Value = "help@bilanceapp.com"
$ws.Range("E3").Value = "fintech"
$ws.Range("F3").Value = "https://play.google.com/store/apps/details?id=com.bilanceapp.flutterapp"
$ws.Range("G3").Value = "https://bilanceapp.com/"

$ws.Range("A4").Value = "Allio Finance"
$ws.Range("B4").Value = "Invest your change, change your life — with Allio.`nAllio is the only finance app powered by machines, and managed by experts to help you automatically save money and build wealth in the 21st Century. Get started today with as little as `$10. `n8 REASONS YOU’LL LOVE ALLIO`n1. Getting started is quick and painless. You’ll be on your way to creating a brighter financial future in under 5 minutes.`n2. Saving & investing is easy and stress-free (because Allio is 100% automated). Set your account up with just a few clicks, and we’ll handle the rest. `n3. You can start building wealth now (even if you don't have much money to spare). With Allio, you can start investing with as little as `$10.`n4. Have trouble saving money? With Allio, your everyday purchases can help you reach your financial goals. With spare change round-ups and multipliers, we’ll help you save without having to think about it.`n5. You’ll be investing like the 1%. No matter how much you have to invest, you’ll get access to Allio’s automated investment portfolios. Crafted by veteran finance experts, these global macro portfolios are highly diversified to reduce risk and increase the potential for long-term growth.`n6. We’ll help you reach (all) your financial goals. Want to save for your rainy day fund, a new car, your next vacation, and retirement? No problem. With Allio, goals are unlimited. We’ll build your portfolio(s) based on when you need the money, and your risk tolerance. `n7. You’ll stay true to your values while maximizing financial returns. With Allio, you can invest in your values, including Clean Energy, Cancer Research, Social Justice & Equity, Made in America, and more.`n8. Never wonder how much money you have again. Track your net worth in real-time to see the progress you're making toward your goals.`nAllio makes saving and investing easy for anyone, and accessible to everyone. With a few simple steps, you can start investing in the future you want. Why wait? Start investing today and take control of your financial future.`nSave. Invest. Grow. With Allio."
$ws.Range("C4").Value = "50,000+"
$ws.Range("D4").Value = "support@alliofinance.com"
$ws.Range("E4").Value = "fintech"
$ws.Range("F4").Value = "https://play.google.com/store/apps/details?id=com.alliofinance.app"
$ws.Range("G4").Value = "https://www.alliofinance.com/"

$ws.Range("A5").Value = "Skip - Fund Your Business"
$ws.Range("B5").Value = "Go further and unlock your full entrepreneurial funding potential with Skip! Skip is here to help you successfully discover, apply for, and maximize your chances of getting your business funded. Here's what you can expect in Skip's new and improved app for small businesses and entrepreneurs:`n— Find the latest funding opportunities, including grants, loans, and other financing options to help you grow your business, whether you're just starting out or you already have a successful business `n— Turn on Skip's grants auto apply to have Skip automatically apply for grants on your behalf every single month`n— Ask all your business growth and funding-related questions on our new Skip AI`n— Maximize your chances of receiving funding opportunities by tracking funding on your Skip Dashboard and writing your applications with Skip AI. You'll get automatic feedback and tips on grant applications. Plus you can request 1-on-1 assistance from Skip experts.`n— Watch our latest videos on small business and funding help, including access to app-only and exclusive content`n— Stay in the know on all things small business through our latest blog posts`n— Join the fastest-growing community supporting entrepreneurs and small businesses in the US"
$ws.Range("C5").Value = "50,000+"
$ws.Range("D5").Value = "hello@helloskip.com"
$ws.Range("E5").Value = "mobile dashboards"
$ws.Range("F5").Value = "https://play.google.com/store/apps/details?id=com.helloskip"
$ws.Range("G5").Value = "http://helloskip.com"

$ws.Range("A6").Value = "Business-in-a-Box"
$ws.Range("B6").Value = "All the templates you need to Start, Plan, Organize, Manage, Finance & Grow your business, in one place.`nBusiness in a Box - The #1 Document Management Software!`nQuickly and easily create any business or legal document with over 2,600 templates.`nBusiness in a Box helps businesspeople get more done in less time. This do-it-yourself document templates software has been created to increase your productivity and efficiency, so you can invest the time and money you save into your business. Today, Business in a Box is used by more than 12 million entrepreneurs, business owners, managers, and executives in over 195 countries.`nBusiness in a Box includes:`nBusiness plan kit, financial statements, business policies and procedures, HR documents, contracts and agreements, board resolutions, letters and notices, checklists, forms ... and more.`nGet organized with the Business in a Box Cloud Drive. Bring all your important files and folders together in one place.`nCustomize templates with the Business in a Box Editor, or use any other document editing software like Microsoft Word, Google Docs, or Pages.`nEdit, store, share and collaborate on all your business files and folders."
$ws.Range("C6").Value = "10,000+"
$ws.Range("D6").Value = "support@biztree.com"
$ws.Range("E6").Value = "mobile dashboards"
$ws.Range("F6").Value = "https://play.google.com/store/apps/details?id=com.biztree.Business_in_a_Box"
$ws.Range("G6").Value = "https://www.business-in-a-box.com"

$ws.Range("A7").Value = "Centro Business Planning Tool"
$ws.Range("B7").Value = "Have fun creating your business plan through simple activities and easily discover and connect with lenders and resources in your area!`nTry our guided step-by-step process to create your basic business plan and budget. We will organize your thinking and email you your first methodical and complete business plan with financials so you can prepare to launch or grow your business. You will be able to track your progress thanks to the To-Do List integrated into the app which will help you understand the steps ahead.`nUnique activities that focus on:`n- Mission, Vision, Values`n- Personal Finance`n- Market Analysis`n- Products and Services`n- Marketing Strategy`n- Operations`n- Business Finances`nThe Centro App also helps you find the best micro-lenders for your business based on your location, your business stage, how much you need, and your industry experience. We’ll connect you with providers that work with small business owners in the United States, and are nonprofit lenders. `nFrom marketing consulting to credit health advisors to legal aid, the app also connects you to additional business resources based on where you are and what you need. The Centro App connects you with providers that work with small business owners in the United States, and speak multiple languages.`nLastly, the To-Do List, where you can track your progress and understand the steps ahead. Each section of the Business Plan you are creating has its own To-Do List, that will make it easier to understand the tasks that you need to complete for your business."
$ws.Range("C7").Value = "50,000+"
$ws.Range("D7").Value = "app@centrocommunity.org"
$ws.Range("E7").Value = "mobile dashboards"
$ws.Range("F7").Value = "https://play.google.com/store/apps/details?id=centrocommunity.org"
$ws.Range("G7").Value = "http://centro-mobile.org/"

